$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A3").Value = "Valkyrie's Sword"
